$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") values repulled/updated for these rows.
# Row => new value
$updates = @{
    3  = 1
    4  = -5
    5  = 2
    6  = -1
    7  = -2
    8  = -1
    9  = 1
    10 = 2
    11 = 8
    13 = 2
    14 = -1
    16 = -4
    17 = 2
    19 = -1
    20 = -4
    22 = -1
    23 = 6
    24 = 6
    25 = -4
    26 = 4
    28 = 3
    29 = -4
    30 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
